$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 95.333336
$ws.Range("J2").Value = 80
$ws.Range("L2").Value = 80
$ws.Range("N2").Value = -306

$ws.Range("H4").Value = 3145.7144
$ws.Range("I4").Value = 3145.7144
$ws.Range("K4").Value = 3145.7144
$ws.Range("M4").Value = -3031.7144

$ws.Range("H5").Value = 103
$ws.Range("I5").Value = 87.8125
$ws.Range("K5").Value = 87.8125
$ws.Range("M5").Value = 27.1875

$ws.Range("H17").Value = 1630.1428
$ws.Range("J17").Value = 2619.6667
$ws.Range("L17").Value = 7859.000100000001
$ws.Range("N17").Value = -8195.000100000001

$ws.Range("H40").Value = 7248.1665
$ws.Range("J40").Value = 7248.1665
$ws.Range("L40").Value = 7248.1665
$ws.Range("N40").Value = -7598.1665

$ws.Range("H54").Value = 18333.334
$ws.Range("I54").Value = 10000
$ws.Range("K54").Value = 10000
$ws.Range("M54").Value = -9514

$ws.Range("H80").Value = 631.4286
$ws.Range("I80").Value = 595
$ws.Range("J80").Value = 646
$ws.Range("K80").Value = 1785
$ws.Range("L80").Value = 1938
$ws.Range("M80").Value = -787
$ws.Range("N80").Value = -3934

$ws.Range("H83").Value = 631.4286
$ws.Range("I83").Value = 595
$ws.Range("J83").Value = 646
$ws.Range("K83").Value = 5355
$ws.Range("L83").Value = 5814
$ws.Range("M83").Value = -363
$ws.Range("N83").Value = -15798

$ws.Range("H86").Value = 9499.5
$ws.Range("I86").Value = 8999
$ws.Range("J86").Value = 10000
$ws.Range("K86").Value = 8999
$ws.Range("L86").Value = 10000
$ws.Range("M86").Value = -7876
$ws.Range("N86").Value = -12246

$ws.Range("H89").Value = 9499.5
$ws.Range("I89").Value = 8999
$ws.Range("J89").Value = 10000
$ws.Range("K89").Value = 44995
$ws.Range("L89").Value = 50000
$ws.Range("M89").Value = -39379
$ws.Range("N89").Value = -61232

$ws.Range("H129").Value = 1319.25
$ws.Range("J129").Value = 0
$ws.Range("L129").Value = 0
$ws.Range("N129").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 912.6
$ws.Range("I2").Value = 975.1111
$ws.Range("K2").Value = 975.1111
$ws.Range("M2").Value = -862.1111

$ws.Range("H4").Value = 283.66666
$ws.Range("I4").Value = 283.66666
$ws.Range("K4").Value = 283.66666
$ws.Range("M4").Value = -167.66666

$ws.Range("H39").Value = 0
$ws.Range("I39").Value = 0
$ws.Range("K39").Value = 0
$ws.Range("M39").ClearContents()

$ws.Range("H102").Value = 783.75
$ws.Range("I102").Value = 781.4286
$ws.Range("J102").Value = 800
$ws.Range("K102").Value = 781.4286
$ws.Range("L102").Value = 800
$ws.Range("M102").Value = 840.5714
$ws.Range("N102").Value = -4044

$ws.Range("H116").Value = 912.6
$ws.Range("I116").Value = 975.1111
$ws.Range("K116").Value = 975.1111
$ws.Range("M116").Value = 1318.8889

$ws.Range("H132").Value = 1537.25
$ws.Range("J132").Value = 1599.5
$ws.Range("L132").Value = 4798.5
$ws.Range("N132").Value = -9858.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 912.6
$ws.Range("I3").Value = 975.1111
$ws.Range("K3").Value = 975.1111
$ws.Range("M3").Value = -861.1111

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 280.7647
$ws.Range("I7").Value = 325.72726
$ws.Range("K7").Value = 325.72726
$ws.Range("M7").Value = -212.72726

$ws.Range("H14").Value = 425.5
$ws.Range("I14").Value = 550
$ws.Range("J14").Value = 301
$ws.Range("K14").Value = 550
$ws.Range("L14").Value = 301
$ws.Range("M14").Value = -380
$ws.Range("N14").Value = -641

$ws.Range("H86").Value = 3014.5
$ws.Range("I86").Value = 2817.4
$ws.Range("J86").Value = 4000
$ws.Range("K86").Value = 2817.4
$ws.Range("L86").Value = 4000
$ws.Range("M86").Value = -1694.4
$ws.Range("N86").Value = -6246

$ws.Range("H89").Value = 3014.5
$ws.Range("I89").Value = 2817.4
$ws.Range("J89").Value = 4000
$ws.Range("K89").Value = 14087
$ws.Range("L89").Value = 20000
$ws.Range("M89").Value = -8471
$ws.Range("N89").Value = -31232

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 143019
$ws.Range("I4").Value = 123.5
$ws.Range("K4").Value = 370.5
$ws.Range("M4").Value = -258.5

$ws.Range("H59").Value = 300
$ws.Range("I59").Value = 300
$ws.Range("K59").Value = 900
$ws.Range("M59").Value = -360

$ws.Range("H121").Value = 1820.2858
$ws.Range("I121").Value = 602.25
$ws.Range("K121").Value = 1806.75
$ws.Range("M121").Value = -496.75

$ws.Range("H122").Value = 5171.5454
$ws.Range("J122").Value = 4652.1665
$ws.Range("L122").Value = 41869.4985
$ws.Range("N122").Value = -46769.4985

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H14").Value = 3572828.5
$ws.Range("I14").Value = 7500900
$ws.Range("J14").Value = 2001600
$ws.Range("K14").Value = 7500900
$ws.Range("L14").Value = 2001600
$ws.Range("M14").Value = -7500732
$ws.Range("N14").Value = -2001936

$ws.Range("H44").Value = 0
$ws.Range("J44").Value = 0
$ws.Range("L44").Value = 0
$ws.Range("N44").ClearContents()

$ws.Range("H46").Value = 18750
$ws.Range("I46").Value = 15000
$ws.Range("K46").Value = 15000
$ws.Range("M46").Value = -14844

$ws.Range("H80").Value = 2423.8333
$ws.Range("I80").Value = 2308.6
$ws.Range("K80").Value = 2308.6
$ws.Range("M80").Value = -1310.6

$ws.Range("H83").Value = 2423.8333
$ws.Range("I83").Value = 2308.6
$ws.Range("K83").Value = 11543
$ws.Range("M83").Value = -6551

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 27812.375
$ws.Range("I132").Value = 32999.8
$ws.Range("J132").Value = 19166.666
$ws.Range("K132").Value = 98999.40000000001
$ws.Range("L132").Value = 57499.99800000001
$ws.Range("M132").Value = -96469.40000000001
$ws.Range("N132").Value = -62559.99800000001

$ws.Range("H136").Value = 7150.1665
$ws.Range("I136").Value = 3875.25
$ws.Range("J136").Value = 13700
$ws.Range("K136").Value = 11625.75
$ws.Range("L136").Value = 41100
$ws.Range("M136").Value = -9075.75
$ws.Range("N136").Value = -46200

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H4").Value = 1070
$ws.Range("I4").Value = 93.333336
$ws.Range("J4").Value = 4000
$ws.Range("K4").Value = 93.333336
$ws.Range("L4").Value = 4000
$ws.Range("M4").Value = 19.666664
$ws.Range("N4").Value = -4226

$ws.Range("H81").Value = 7640.1665
$ws.Range("I81").Value = 8043.8184
$ws.Range("K81").Value = 16087.6368
$ws.Range("M81").Value = -15026.6368

$ws.Range("H84").Value = 7640.1665
$ws.Range("I84").Value = 8043.8184
$ws.Range("K84").Value = 80438.18400000001
$ws.Range("M84").Value = -75134.18400000001
